$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before column D (shifts D:K -> F:M)
$ws.Range("D:E").Insert()

# Copy number formats from the (now-shifted) old columns F:M into new D:E
$ws.Range("F7:M7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F8:M102").Copy()
$ws.Range("D8:E102").PasteSpecial(-4122)
$ws.Range("F38:M38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F80:M80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new D:E columns with the newest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 4500
$ws.Range("E8").Value = 5300
$ws.Range("D9").Value = 1100
$ws.Range("E9").Value = 1400
$ws.Range("D10").Value = 3400
$ws.Range("E10").Value = 3900
$ws.Range("D12").Value = 7900
$ws.Range("E12").Value = 7200
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = -109800
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 1700
$ws.Range("E15").Value = 1900
$ws.Range("D17").Value = -82900
$ws.Range("E17").Value = 30100
$ws.Range("D18").Value = 87400
$ws.Range("E18").Value = -24800
$ws.Range("D20").Value = 200
$ws.Range("E20").Value = 7600
$ws.Range("D21").Value = 89200
$ws.Range("E21").Value = -15300
$ws.Range("D22").Value = 3000
$ws.Range("E22").Value = 3400
$ws.Range("D23").Value = 84600
$ws.Range("E23").Value = -20600
$ws.Range("D24").Value = 500
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 84100
$ws.Range("E26").Value = -20600
$ws.Range("D27").Value = 84100
$ws.Range("E27").Value = -20600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -200
$ws.Range("E32").Value = -7600
$ws.Range("D33").Value = 84100
$ws.Range("E33").Value = -20600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 84100
$ws.Range("E35").Value = -20600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 122500
$ws.Range("E41").Value = 67400
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 1600
$ws.Range("E43").Value = 2500
$ws.Range("D44").Value = 3900
$ws.Range("E44").Value = 6300
$ws.Range("D45").Value = 4200
$ws.Range("E45").Value = 2300
$ws.Range("D46").Value = 132300
$ws.Range("E46").Value = 78500
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 300
$ws.Range("E48").Value = 300
$ws.Range("D49").Value = 37400
$ws.Range("E49").Value = 61700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 300
$ws.Range("E52").Value = 300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 170300
$ws.Range("E54").Value = 140800
$ws.Range("D57").Value = 1200
$ws.Range("E57").Value = 3200
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 16100
$ws.Range("E59").Value = 20000
$ws.Range("D60").Value = 17200
$ws.Range("E60").Value = 23200
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 79100
$ws.Range("D62").Value = 40100
$ws.Range("E62").Value = 44800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 57300
$ws.Range("E66").Value = 147100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -211000
$ws.Range("E72").Value = -295100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 113000
$ws.Range("E76").Value = -6300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 84100
$ws.Range("E81").Value = -20600
$ws.Range("D83").Value = 1700
$ws.Range("E83").Value = 1900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -32700
$ws.Range("E89").Value = -19300
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 159300
$ws.Range("E94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -71500
$ws.Range("E100").Value = 1200
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 55100
$ws.Range("E102").Value = -18100

# Row 91 (Capital Expenditures): mark older quarters without data as "NA"
$ws.Range("F91").Value = "NA"
$ws.Range("G91").Value = "NA"
$ws.Range("H91").Value = "NA"
$ws.Range("I91").Value = "NA"
$ws.Range("J91").Value = "NA"
